# "replaced na with NA"
# Every cell that used to display the shared string "na" is updated to
# display "NA" instead -- except cell S38, which keeps the original "na".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TrayPlaybackExperimentData")

$cellsToFix = @(
    "P2","P3","L4","P4","L5","L8","L11","L14","L20","L29","P29",
    "L30","L31","L32","L33","L34","L53","L54","L56","L60","L63",
    "P69","P70","P71","P79"
)

foreach ($addr in $cellsToFix) {
    $ws.Range($addr).Value = "NA"
}

# Reflect the updated selection state seen in the saved workbook.
$ws.Range("P:P").Select()
